# Autogenerated on Thu Mar 26 2015 18:06:15 GMT+0000 (Coordinated Universal Time)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the hyperlink that used to sit on A72 (the Eurostat database URL).
$ws.Range("A72").Hyperlinks.Delete()

# --- First "Source:" block (rows 70-74) ---
# Row 70 "Source:" is unchanged.

# Row 71 becomes a blank spacer line (it used to hold the "SBS Main
# Indicators..." text, which now moves down to row 72). It already has the
# italic "source" look, so just drop its text.
$ws.Range("A71").ClearContents()

# Row 72 now holds the "SBS Main Indicators..." text that used to be on row
# 71. The cell used to be the hyperlink ("HyperLink" style: underlined,
# blue) - switch it to the plain italic "source" look instead.
$ws.Range("A72").Value2 = "SBS Main Indicators, Annual enterprise statistics by size class for special aggregates of activities (NACE Rev. 2)"
$ws.Range("A72").Font.Underline = -4142
$ws.Range("A72").Font.ColorIndex = -4105
$ws.Range("A72").Font.Bold = $false
$ws.Range("A72").Font.Italic = $true

# Row 73 stays the blank spacer line that separates the text from the URL -
# already correct, leave untouched.

# Row 74 is a new row holding the plain-text URL (no hyperlink any more).
$ws.Range("A74").Value2 = "http://epp.eurostat.ec.europa.eu/portal/page/portal/european_business/data/database"
$ws.Range("A74").Font.Italic = $true

# --- Second source block (rows 77-80, used to be rows 76-79) ---
# Old row 76 ("Statistical Office of the Republic of Slovenia", title style)
# moves down to row 77, so drop the old row 76 cell entirely.
$ws.Range("A76").Clear()

# Row 77: title "Statistical Office of the Republic of Slovenia".
$ws.Range("A77").Value2 = "Statistical Office of the Republic of Slovenia"
$ws.Range("A77").Font.Bold = $true
$ws.Range("A77").Font.Italic = $false

# Row 78: the same text repeated using the italic "source" look, replacing
# the old (incorrect) long SI-STAT description.
$ws.Range("A78").Value2 = "Statistical Office of the Republic of Slovenia"
$ws.Range("A78").Font.Bold = $false
$ws.Range("A78").Font.Italic = $true

# Row 79: title "SBS Eurostat", replacing the old "SI-STAT Data Portal..."
# text that used to live here.
$ws.Range("A79").Value2 = "SBS Eurostat"
$ws.Range("A79").Font.Bold = $true
$ws.Range("A79").Font.Italic = $false

# Row 80 is a new row: "SBS Eurostat" repeated using the italic "source"
# look, replacing the old (incorrect) "Sructural Business Statistics..."
# description.
$ws.Range("A80").Value2 = "SBS Eurostat"
$ws.Range("A80").Font.Bold = $false
$ws.Range("A80").Font.Italic = $true
